$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to Text so they
# are stored as strings (matching the scraped-data formatting), not coerced
# into Excel numbers/percentages.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "259.47"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.66%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "20"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.95"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.52%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "20"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.702"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.60%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "20"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06032"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.50%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "20"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.674"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.51%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "20"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8585"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.10%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "20"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9297"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.33%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "20"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1398"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.70%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "20"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04756"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "20.38%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "20"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07055"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.27%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "20"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03121"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.74%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "20"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09140"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.29%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "20"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001528"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.22%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "20"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006035"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.18%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "20"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006072"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.42%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "20"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.466"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "20"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.156"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-1.47%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.67%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "20"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.26%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "20"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1296"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.36%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "20"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.129"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.25%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "20"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04231"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.69%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "20"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.30%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "20"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004051"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-5.73%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "20"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.17%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "20"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001520"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-21.46%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "20"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "20"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "20"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "20"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "20"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "20"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "20"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "20"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "20"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "20"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "20"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "20"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "20"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03856"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.21%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "20"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.28%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "20"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.003928"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.81%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "20"
$ws.Range("B43").Value = "LocalTraders"
$ws.Range("C43").Value = "https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.01527"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "29.77%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "20"
$ws.Range("B44").Value = "CEJI"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002290"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.77%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "20"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005086"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-6.68%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "20"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.14%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "20"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1303"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-4.09%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "20"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05443"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-9.20%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "20"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.14%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "20"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.14%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "20"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "20"
